# Onboarding3FnStep: add "OnBoardingStep7", "Sheet7" and "OnBoarding" test-data sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new worksheets at the right tab positions:
#    Signup, OnBoardingStep7, Sheet7, LoginData, Dashboard, EditYourProfile, OnBoarding
# ---------------------------------------------------------------------------
$signup = $wb.Worksheets.Item("Signup")
$onboardingStep7 = $wb.Worksheets.Add($null, $signup)
$onboardingStep7.Name = "OnBoardingStep7"

$sheet7 = $wb.Worksheets.Add($null, $onboardingStep7)
$sheet7.Name = "Sheet7"

$editYourProfile = $wb.Worksheets.Item("EditYourProfile")
$onboarding = $wb.Worksheets.Add($null, $editYourProfile)
$onboarding.Name = "OnBoarding"

# ---------------------------------------------------------------------------
# 2. OnBoardingStep7 - weight / height / date validation test data
# ---------------------------------------------------------------------------
$ws = $onboardingStep7

$ws.Range("A1").Value = "validweight"
$ws.Range("B1").Value = "validheight"
$ws.Range("C1").Value = "Invalidweight"
$ws.Range("D1").Value = "Invalidheight"
$ws.Range("E1").Value = "validDate"
$ws.Range("F1").Value = "InvalidDate"
$ws.Range("A1:D1").Font.Name = "Arial"
$ws.Range("E1:F1").Font.Name = "Calibri"

$ws.Range("A2").Value = 34
$ws.Range("B2").Value = 170
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = -4
$ws.Range("A2:D2").Font.Name = "Arial"
$ws.Range("A2:D2").HorizontalAlignment = -4152

# ---------------------------------------------------------------------------
# 3. Sheet7 stays empty (placeholder sheet created by the export)
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 4. OnBoarding - full onboarding form test data
# ---------------------------------------------------------------------------
$ws = $onboarding

$headers = @("UserName","Password","With Report/Without Report","Health Conditions","First Name","Age","BP Status","Menstrual Cycle Track","Last Period Date","Weight in KG","Height in CM","Dietary Preferences`n","Physical Activity Level","Food Allergies`n& Sensitivities","Please select any foods you need to avoid:`n","Medications & Supplements`n","select or add the medications/supplements you're taking","Plan")

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Range($cols[$i] + "1")
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.Font.Name = "Calibri"
    $cell.Interior.Color = 16308937
    $cell.Borders.LineStyle = 1
    $cell.Borders.Color = 0
    $cell.WrapText = $true
}

$extraCols = @("S","T","U","V","W","X","Y","Z")
foreach ($c in $extraCols) {
    $cell = $ws.Range($c + "1")
    $cell.Font.Name = "Calibri"
    $cell.Interior.Color = 16308937
}

# Row 2 - User1
$ws.Range("A2").Value = "User1@gmail.com"
$ws.Range("B2").Value = "User1@987&"
$ws.Range("C2").Value = "without report"
$ws.Range("D2").Value = "Pre-diabetes / Diabetes , Hypothyroidism"
$ws.Range("E2").Value = "User1"
$ws.Range("F2").Value = 43
$ws.Range("G2").Value = "I have been diagnosed but don't take medication"
$ws.Range("H2").Value = "Yes"
$ws.Range("I2").Value = 46031
$ws.Range("J2").Value = 67
$ws.Range("K2").Value = 145
$ws.Range("L2").Value = "Non-Vegetarian"
$ws.Range("M2").Value = "Lightly Active"
$ws.Range("N2").Value = "Yes, I have food intolerances/sensitivities"
$ws.Range("O2").Value = "Fish, Peanuts"
$ws.Range("P2").Value = "Yes, I take prescription medications"
$ws.Range("Q2").Value = "Thyroid medication,`nVitamin D, Insulin"
$ws.Range("R2").Value = "Free"

# Row 3 - User2
$ws.Range("A3").Value = "User2@gmail.com"
$ws.Range("B3").Value = "User2*/245"
$ws.Range("C3").Value = "without report"
$ws.Range("D3").Value = "PCOS"
$ws.Range("E3").Value = "User2"
$ws.Range("F3").Value = 25
$ws.Range("G3").Value = "I'm not sure"
$ws.Range("H3").Value = "My cycle is irregular"
$ws.Range("I3").Value = 46054
$ws.Range("J3").Value = 56
$ws.Range("K3").Value = 160
$ws.Range("L3").Value = "Vegetarian Diet"
$ws.Range("M3").Value = "Sedentary"
$ws.Range("N3").Value = "No, I can eat everything"
$ws.Range("P3").Value = "No, I don't take any medications or supplements"
$ws.Range("R3").Value = "Free"

# Formatting for the two data rows (A:R) + blank tail columns (S:Z)
$dataCols = @("A","B","C","D","E","G","H","L","M","N","O","P","Q","R")
$numCols = @("F","J","K")
foreach ($r in 2,3) {
    foreach ($c in $dataCols) {
        $cell = $ws.Range($c + $r)
        $cell.Font.Name = "Calibri"
        $cell.Borders.LineStyle = 1
        $cell.Borders.Color = 0
        $cell.WrapText = $true
    }
    foreach ($c in $numCols) {
        $cell = $ws.Range($c + $r)
        $cell.Font.Name = "Calibri"
        $cell.Borders.LineStyle = 1
        $cell.Borders.Color = 0
        $cell.WrapText = $true
        $cell.HorizontalAlignment = -4152
    }
    $dateCell = $ws.Range("I" + $r)
    $dateCell.Font.Name = "Calibri"
    $dateCell.Borders.LineStyle = 1
    $dateCell.Borders.Color = 0
    $dateCell.WrapText = $true
    $dateCell.HorizontalAlignment = -4152
    $dateCell.NumberFormat = "mm/dd/yyyy"

    foreach ($c in $extraCols) {
        $cell = $ws.Range($c + $r)
        $cell.Font.Name = "Calibri"
    }
}

# The two intentionally blank-but-bordered cells on row 3 (O3, Q3)
$ws.Range("O3").Font.Name = "Calibri"
$ws.Range("O3").Borders.LineStyle = 1
$ws.Range("O3").Borders.Color = 0
$ws.Range("Q3").Font.Name = "Calibri"
$ws.Range("Q3").Borders.LineStyle = 1
$ws.Range("Q3").Borders.Color = 0

Write-Host "Onboarding sheets added."
